# egresos.xlsx fix: "est consideraba 300 peores"
# Update projected values in columns B (Curso), C (Futuras), D (SEM) and E (Total)
# for rows 2-95 with corrected actuarial estimates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1875064470.3247409
$ws.Range("C2").Value = 273024684.85272998
$ws.Range("D2").Value = 141462564.25284281
$ws.Range("E2").Value = 2289551719.4303141
$ws.Range("B3").Value = 1764732611.3479791
$ws.Range("C3").Value = 497199341.04982382
$ws.Range("D3").Value = 142142877.65777829
$ws.Range("E3").Value = 2404074830.0555811
$ws.Range("B4").Value = 1656591368.474839
$ws.Range("C4").Value = 775837193.26067078
$ws.Range("D4").Value = 146168230.77330899
$ws.Range("E4").Value = 2578596792.5088191
$ws.Range("B5").Value = 1554963587.0574341
$ws.Range("C5").Value = 1147089902.9605441
$ws.Range("D5").Value = 156982485.71952859
$ws.Range("E5").Value = 2859035975.7375059
$ws.Range("B6").Value = 1462868885.2175331
$ws.Range("C6").Value = 1525798925.3698151
$ws.Range("D6").Value = 175113842.71487769
$ws.Range("E6").Value = 3163781653.3022261
$ws.Range("B7").Value = 1369160247.033267
$ws.Range("C7").Value = 1820613777.9863341
$ws.Range("D7").Value = 188565056.59482411
$ws.Range("E7").Value = 3378339081.6144252
$ws.Range("B8").Value = 1282903565.055856
$ws.Range("C8").Value = 2214269483.2385778
$ws.Range("D8").Value = 200510226.36540321
$ws.Range("E8").Value = 3697683274.6598382
$ws.Range("B9").Value = 1199594784.064779
$ws.Range("C9").Value = 2540979080.456326
$ws.Range("D9").Value = 225958006.18789169
$ws.Range("E9").Value = 3966531870.7089958
$ws.Range("B10").Value = 1117146585.3497331
$ws.Range("C10").Value = 2775828073.2979031
$ws.Range("D10").Value = 240732496.10213101
$ws.Range("E10").Value = 4133707154.7497659
$ws.Range("B11").Value = 1047352320.7850699
$ws.Range("C11").Value = 3095257935.0111461
$ws.Range("D11").Value = 259046337.0454123
$ws.Range("E11").Value = 4401656592.8416281
$ws.Range("B12").Value = 975111042.19665504
$ws.Range("C12").Value = 3406542587.8199492
$ws.Range("D12").Value = 278630545.48306942
$ws.Range("E12").Value = 4660284175.4996738
$ws.Range("B13").Value = 908610399.51434541
$ws.Range("C13").Value = 3680422146.8360929
$ws.Range("D13").Value = 297963155.00074452
$ws.Range("E13").Value = 4886995701.3511839
$ws.Range("B14").Value = 846465984.12709343
$ws.Range("C14").Value = 3987617368.8154688
$ws.Range("D14").Value = 317054856.75896531
$ws.Range("E14").Value = 5151138209.7015276
$ws.Range("B15").Value = 786784130.46574128
$ws.Range("C15").Value = 4249982907.8708401
$ws.Range("D15").Value = 336018528.61071831
$ws.Range("E15").Value = 5372785566.9473
$ws.Range("B16").Value = 732096322.34113419
$ws.Range("C16").Value = 4427005240.5263929
$ws.Range("D16").Value = 352103573.25092173
$ws.Range("E16").Value = 5511205136.1184483
$ws.Range("B17").Value = 679417834.48389888
$ws.Range("C17").Value = 4700378000.384244
$ws.Range("D17").Value = 376712840.36626017
$ws.Range("E17").Value = 5756508675.2344036
$ws.Range("B18").Value = 629113315.07975876
$ws.Range("C18").Value = 4947207926.8931046
$ws.Range("D18").Value = 392995119.49281669
$ws.Range("E18").Value = 5969316361.4656801
$ws.Range("B19").Value = 580912503.96303701
$ws.Range("C19").Value = 5182069393.8398085
$ws.Range("D19").Value = 406441696.54067647
$ws.Range("E19").Value = 6169423594.3435211
$ws.Range("B20").Value = 535339939.97448701
$ws.Range("C20").Value = 5365394255.7475405
$ws.Range("D20").Value = 422237550.01942742
$ws.Range("E20").Value = 6322971745.7414541
$ws.Range("B21").Value = 491784441.3606832
$ws.Range("C21").Value = 5680197937.1307716
$ws.Range("D21").Value = 446133012.40665472
$ws.Range("E21").Value = 6618115390.8981094
$ws.Range("B22").Value = 450213033.30866992
$ws.Range("C22").Value = 5917852971.6784496
$ws.Range("D22").Value = 467985064.90809208
$ws.Range("E22").Value = 6836051069.8952122
$ws.Range("B23").Value = 410613700.70397413
$ws.Range("C23").Value = 6080701651.1342382
$ws.Range("D23").Value = 485715141.58324963
$ws.Range("E23").Value = 6977030493.4214621
$ws.Range("B24").Value = 372959422.56024587
$ws.Range("C24").Value = 6212220206.4108133
$ws.Range("D24").Value = 500416140.29246378
$ws.Range("E24").Value = 7085595769.2635231
$ws.Range("B25").Value = 337215274.18982953
$ws.Range("C25").Value = 6299364322.788866
$ws.Range("D25").Value = 513350145.80405009
$ws.Range("E25").Value = 7149929742.7827463
$ws.Range("B26").Value = 303376935.71308851
$ws.Range("C26").Value = 6347749984.2452288
$ws.Range("D26").Value = 515260683.38966322
$ws.Range("E26").Value = 7166387603.3479795
$ws.Range("B27").Value = 271452610.35055411
$ws.Range("C27").Value = 6397731534.6682615
$ws.Range("D27").Value = 519371032.48157543
$ws.Range("E27").Value = 7188555177.500392
$ws.Range("B28").Value = 241458065.20978051
$ws.Range("C28").Value = 6296230510.1168842
$ws.Range("D28").Value = 518623684.97951972
$ws.Range("E28").Value = 7056312260.3061848
$ws.Range("B29").Value = 213419229.50000361
$ws.Range("C29").Value = 6193771509.4332991
$ws.Range("D29").Value = 514199782.07821202
$ws.Range("E29").Value = 6921390521.0115147
$ws.Range("B30").Value = 187356017.14784089
$ws.Range("C30").Value = 6056950094.132823
$ws.Range("D30").Value = 502490499.25048631
$ws.Range("E30").Value = 6746796610.5311499
$ws.Range("B31").Value = 163284504.2168144
$ws.Range("C31").Value = 5911127866.0218678
$ws.Range("D31").Value = 488626881.06197572
$ws.Range("E31").Value = 6563039251.3006573
$ws.Range("B32").Value = 141224866.18251839
$ws.Range("C32").Value = 5721845884.5616379
$ws.Range("D32").Value = 472439181.80986679
$ws.Range("E32").Value = 6335509932.5540228
$ws.Range("B33").Value = 121193652.3915122
$ws.Range("C33").Value = 5482965633.305419
$ws.Range("D33").Value = 452450213.64629883
$ws.Range("E33").Value = 6056609499.3432293
$ws.Range("B34").Value = 103179586.957252
$ws.Range("C34").Value = 5244279418.5521984
$ws.Range("D34").Value = 430717603.16950673
$ws.Range("E34").Value = 5778176608.678957
$ws.Range("B35").Value = 87146791.089476481
$ws.Range("C35").Value = 4981733204.3188868
$ws.Range("D35").Value = 408326980.70833111
$ws.Range("E35").Value = 5477206976.1166945
$ws.Range("B36").Value = 73042001.80689466
$ws.Range("C36").Value = 4722044288.8066959
$ws.Range("D36").Value = 386258407.45599627
$ws.Range("E36").Value = 5181344698.0695868
$ws.Range("B37").Value = 60787925.899342127
$ws.Range("C37").Value = 4445827039.1752167
$ws.Range("D37").Value = 362982577.05473399
$ws.Range("E37").Value = 4869597542.1292934
$ws.Range("B38").Value = 50278812.433501117
$ws.Range("C38").Value = 4173666906.673368
$ws.Range("D38").Value = 339756641.93762797
$ws.Range("E38").Value = 4563702361.0444965
$ws.Range("B39").Value = 41382571.273546956
$ws.Range("C39").Value = 3906753209.805963
$ws.Range("D39").Value = 317208768.60380572
$ws.Range("E39").Value = 4265344549.6833162
$ws.Range("B40").Value = 33945194.73889415
$ws.Range("C40").Value = 3642298584.3753762
$ws.Range("D40").Value = 295270346.99934667
$ws.Range("E40").Value = 3971514126.113616
$ws.Range("B41").Value = 27794318.761961222
$ws.Range("C41").Value = 3389364415.4506078
$ws.Range("D41").Value = 273476201.08168209
$ws.Range("E41").Value = 3690634935.294251
$ws.Range("B42").Value = 22758324.692854811
$ws.Range("C42").Value = 3137402962.9047298
$ws.Range("D42").Value = 252220237.22631091
$ws.Range("E42").Value = 3412381524.823895
$ws.Range("B43").Value = 18665328.684356909
$ws.Range("C43").Value = 2895785980.631938
$ws.Range("D43").Value = 232380717.58524209
$ws.Range("E43").Value = 3146832026.9015369
$ws.Range("B44").Value = 15351621.152417749
$ws.Range("C44").Value = 2666168254.1701331
$ws.Range("D44").Value = 213348809.2903595
$ws.Range("E44").Value = 2894868684.6129098
$ws.Range("B45").Value = 12669520.57881755
$ws.Range("C45").Value = 2445897343.7148242
$ws.Range("D45").Value = 195075249.1387521
$ws.Range("E45").Value = 2653642113.432394
$ws.Range("B46").Value = 10492127.392676281
$ws.Range("C46").Value = 2236978796.6671839
$ws.Range("D46").Value = 177677382.83587179
$ws.Range("E46").Value = 2425148306.8957319
$ws.Range("B47").Value = 8714773.4957074188
$ws.Range("C47").Value = 2039884675.236114
$ws.Range("D47").Value = 161328082.3392722
$ws.Range("E47").Value = 2209927531.071094
$ws.Range("B48").Value = 7253373.5581587721
$ws.Range("C48").Value = 1854411576.4980831
$ws.Range("D48").Value = 146009114.69860169
$ws.Range("E48").Value = 2007674064.754843
$ws.Range("B49").Value = 6042432.6929495912
$ws.Range("C49").Value = 1680160976.886025
$ws.Range("D49").Value = 131678864.4047624
$ws.Range("E49").Value = 1817882273.983737
$ws.Range("B50").Value = 5031742.4097335935
$ws.Range("C50").Value = 1516737557.842576
$ws.Range("D50").Value = 118297442.20769
$ws.Range("E50").Value = 1640066742.4599991
$ws.Range("B51").Value = 4183219.59936509
$ws.Range("C51").Value = 1363814232.6976359
$ws.Range("D51").Value = 105832309.32765
$ws.Range("E51").Value = 1473829761.624651
$ws.Range("B52").Value = 3467929.0914756618
$ws.Range("C52").Value = 1221119322.0808439
$ws.Range("D52").Value = 94257476.307642817
$ws.Range("E52").Value = 1318844727.4799631
$ws.Range("B53").Value = 2863315.3700107788
$ws.Range("C53").Value = 1088367967.7254951
$ws.Range("D53").Value = 83544417.677061319
$ws.Range("E53").Value = 1174775700.772567
$ws.Range("B54").Value = 2351411.6868636529
$ws.Range("C54").Value = 965259697.15346622
$ws.Range("D54").Value = 73663318.312148735
$ws.Range("E54").Value = 1041274427.1524791
$ws.Range("B55").Value = 1917768.5746119269
$ws.Range("C55").Value = 851512163.85427201
$ws.Range("D55").Value = 64586149.380949043
$ws.Range("E55").Value = 918016081.80983293
$ws.Range("B56").Value = 1550848.761228028
$ws.Range("C56").Value = 746889196.55397856
$ws.Range("D56").Value = 56289058.617527381
$ws.Range("E56").Value = 804729103.93273401
$ws.Range("B57").Value = 1241422.5921559981
$ws.Range("C57").Value = 651164700.91972446
$ws.Range("D57").Value = 48748433.916902438
$ws.Range("E57").Value = 701154557.42878294
$ws.Range("B58").Value = 981910.59325146268
$ws.Range("C58").Value = 564072895.85054076
$ws.Range("D58").Value = 41936524.030783869
$ws.Range("E58").Value = 606991330.47457612
$ws.Range("B59").Value = 765986.98358147347
$ws.Range("C59").Value = 485308300.61785781
$ws.Range("D59").Value = 35822392.625649028
$ws.Range("E59").Value = 521896680.22708827
$ws.Range("B60").Value = 588286.62890181597
$ws.Range("C60").Value = 414542976.83493578
$ws.Range("D60").Value = 30372898.042151231
$ws.Range("E60").Value = 445504161.50598878
$ws.Range("B61").Value = 444106.9935460061
$ws.Range("C61").Value = 351438223.16393149
$ws.Range("D61").Value = 25553899.372556109
$ws.Range("E61").Value = 377436229.53003359
$ws.Range("B62").Value = 329133.92652942613
$ws.Range("C62").Value = 295623584.48005968
$ws.Range("D62").Value = 21328872.362181779
$ws.Range("E62").Value = 317281590.76877087
$ws.Range("B63").Value = 239215.55119182781
$ws.Range("C63").Value = 246671762.50970861
$ws.Range("D63").Value = 17656703.452076819
$ws.Range("E63").Value = 264567681.5129773
$ws.Range("B64").Value = 170338.34205522531
$ws.Range("C64").Value = 204106043.6672723
$ws.Range("D64").Value = 14492813.270346491
$ws.Range("E64").Value = 218769195.27967411
$ws.Range("B65").Value = 118711.25504947281
$ws.Range("C65").Value = 167422303.18664649
$ws.Range("D65").Value = 11791321.99200912
$ws.Range("E65").Value = 179332336.43370509
$ws.Range("B66").Value = 80883.488171971912
$ws.Range("C66").Value = 136107293.5088248
$ws.Range("D66").Value = 9506714.1254327372
$ws.Range("E66").Value = 145694891.12242949
$ws.Range("B67").Value = 53813.392919712307
$ws.Range("C67").Value = 109638412.610193
$ws.Range("D67").Value = 7593818.1605135929
$ws.Range("E67").Value = 117286044.1636263
$ws.Range("B68").Value = 34905.00498658049
$ws.Range("C68").Value = 87486071.981144235
$ws.Range("D68").Value = 6007922.0087512974
$ws.Range("E68").Value = 93528898.994882122
$ws.Range("B69").Value = 22034.1725601988
$ws.Range("C69").Value = 69128349.747029021
$ws.Range("D69").Value = 4705960.2680958109
$ws.Range("E69").Value = 73856344.187685028
$ws.Range("B70").Value = 13514.01775488513
$ws.Range("C70").Value = 54068155.5365052
$ws.Range("D70").Value = 3647783.9710556678
$ws.Range("E70").Value = 57729453.525315747
$ws.Range("B71").Value = 8042.5371547073864
$ws.Range("C71").Value = 41843615.996626116
$ws.Range("D71").Value = 2796830.2228980851
$ws.Range("E71").Value = 44648488.756678917
$ws.Range("B72").Value = 4637.9775530884745
$ws.Range("C72").Value = 32028703.88932246
$ws.Range("D72").Value = 2119997.774211328
$ws.Range("E72").Value = 34153339.641086876
$ws.Range("B73").Value = 2587.0217218906018
$ws.Range("C73").Value = 24233775.23544376
$ws.Range("D73").Value = 1587547.7941147459
$ws.Range("E73").Value = 25823910.051280402
$ws.Range("B74").Value = 1390.6029815561351
$ws.Range("C74").Value = 18110178.806836382
$ws.Range("D74").Value = 1173318.269283463
$ws.Range("E74").Value = 19284887.6791014
$ws.Range("B75").Value = 712.34500344354024
$ws.Range("C75").Value = 13353463.622441599
$ws.Range("D75").Value = 854825.32232445513
$ws.Range("E75").Value = 14209001.2897695
$ws.Range("B76").Value = 346.86869514273411
$ws.Range("C76").Value = 9703189.5537421703
$ws.Range("D76").Value = 613091.79540732375
$ws.Range("E76").Value = 10316628.217844641
$ws.Range("B77").Value = 157.88029849020839
$ws.Range("C77").Value = 6938336.3271142328
$ws.Range("D77").Value = 432198.95941688248
$ws.Range("E77").Value = 7370693.1668296056
$ws.Range("C78").Value = 4873832.5439235559
$ws.Range("D78").Value = 298910.15364468843
$ws.Range("E78").Value = 5172809.0512164785
$ws.Range("C79").Value = 3356206.07137352
$ws.Range("D79").Value = 202354.38615965919
$ws.Range("E79").Value = 3558585.9767912431
$ws.Range("C80").Value = 2260131.2915352639
$ws.Range("D80").Value = 133739.64633015401
$ws.Range("E80").Value = 2393879.837626609
$ws.Range("C81").Value = 1484395.995983154
$ws.Range("D81").Value = 86043.683114612533
$ws.Range("E81").Value = 1570442.4732942721
$ws.Range("C82").Value = 947892.69634570309
$ws.Range("D82").Value = 53709.097233808418
$ws.Range("E82").Value = 1001602.57173367
$ws.Range("C83").Value = 586364.4109393541
$ws.Range("D83").Value = 32398.33886368425
$ws.Range("E83").Value = 618762.92611276894
$ws.Range("C84").Value = 349801.90347028471
$ws.Range("D84").Value = 18795.786014294681
$ws.Range("E84").Value = 368597.68948457937
$ws.Range("C85").Value = 200155.82394886229
$ws.Range("D85").Value = 10428.04368764939
$ws.Range("E85").Value = 210583.86763651171
$ws.Range("C86").Value = 109158.9761598325
$ws.Range("D86").Value = 5497.6119644899654
$ws.Range("E86").Value = 114656.5881243225
$ws.Range("C87").Value = 56332.939146573888
$ws.Range("D87").Value = 2734.589907951985
$ws.Range("E87").Value = 59067.529054525883
$ws.Range("C88").Value = 27261.044834391261
$ws.Range("D88").Value = 1272.9494191318779
$ws.Range("E88").Value = 28533.994253523138
$ws.Range("C89").Value = 12235.68768846608
$ws.Range("D89").Value = 549.46577723410019
$ws.Range("E89").Value = 12785.153465700179
$ws.Range("C90").Value = 5031.5164956090948
$ws.Range("D90").Value = 217.83481838858631
$ws.Range("E90").Value = 5249.3513139976812
$ws.Range("C91").Value = 1877.180708975937
$ws.Range("D91").Value = 78.732910104246031
$ws.Range("E91").Value = 1955.913619080183
$ws.Range("C92").Value = 623.23044328852268
$ws.Range("D92").Value = 25.58194999413757
$ws.Range("E92").Value = 648.8123932826602
$ws.Range("C93").Value = 169.3044964436524
$ws.Range("D93").Value = 7.0650618274492478
$ws.Range("E93").Value = 176.36955827110171
$ws.Range("C94").Value = 35.804912239841059
$ws.Range("D94").Value = 1.617431215271647
$ws.Range("E94").Value = 37.422343455112703
$ws.Range("C95").Value = 7.1129140062323506
$ws.Range("D95").Value = 0.34298691981136642
$ws.Range("E95").Value = 7.455900926043717

# Widen columns B:E to match column E's existing best-fit width (12 chars),
# since the new values now share the same order of magnitude as column E.
$ws.Columns("B:D").ColumnWidth = $ws.Columns("E:E").ColumnWidth

# Scroll the view down and clear the stale selection left over from editing,
# matching the post-edit view state (topLeftCell A61, default selection).
$ws.Range("A61").Select() | Out-Null
$ws.Range("A1").Select() | Out-Null
